$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1728506666666667
$ws.Range("H2").Value = 0.518552
$ws.Range("I2").Value = 0.0840503369699626
$ws.Range("J2").Value = 0.0840503369699626
$ws.Range("M2").Value = 0.1352566666666667
$ws.Range("N2").Value = 0.40577
$ws.Range("O2").Value = 0.1173241749329269
$ws.Range("P2").Value = 0.1173241749329268
$ws.Range("Q2").Value = 0.02337920500444444
$ws.Range("R2").Value = 0.21041284504
$ws.Range("S2").Value = 0.009861136437835342
$ws.Range("T2").Value = 0.00986113643783534
$ws.Range("G3").Value = 0.1728506666666667
$ws.Range("H3").Value = 0.518552
$ws.Range("I3").Value = 0.0840503369699626
$ws.Range("J3").Value = 0.0840503369699626
$ws.Range("O3").Value = 0.03951584152489912
$ws.Range("P3").Value = 0.03951584152489912
$ws.Range("Q3").Value = 0.007874327353777778
$ws.Range("R3").Value = 0.07086894618400001
$ws.Range("S3").Value = 0.003321319795819412
$ws.Range("T3").Value = 0.003321319795819412
$ws.Range("G4").Value = 0.1728506666666667
$ws.Range("H4").Value = 0.518552
$ws.Range("I4").Value = 0.0840503369699626
$ws.Range("J4").Value = 0.0840503369699626
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9440163333333332
$ws.Range("N4").Value = 2.832049
$ws.Range("O4").Value = 0.8188575111383802
$ws.Range("P4").Value = 0.8188575111383801
$ws.Range("Q4").Value = 0.1631738525608889
$ws.Range("R4").Value = 1.468564673048
$ws.Range("S4").Value = 0.06882524974156576
$ws.Range("T4").Value = 0.06882524974156576
$ws.Range("G5").Value = 0.1728506666666667
$ws.Range("H5").Value = 0.518552
$ws.Range("I5").Value = 0.0840503369699626
$ws.Range("J5").Value = 0.0840503369699626
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.028017
$ws.Range("N5").Value = 0.084051
$ws.Range("O5").Value = 0.02430247240379386
$ws.Range("P5").Value = 0.02430247240379386
$ws.Range("Q5").Value = 0.004842757128000001
$ws.Range("R5").Value = 0.043584814152
$ws.Range("S5").Value = 0.002042630994742091
$ws.Range("T5").Value = 0.002042630994742091
$ws.Range("I6").Value = 0.6650661694281633
$ws.Range("J6").Value = 0.6650661694281633
$ws.Range("M6").Value = 0.1352566666666667
$ws.Range("N6").Value = 0.40577
$ws.Range("O6").Value = 0.1173241749329269
$ws.Range("P6").Value = 0.1173241749329268
$ws.Range("Q6").Value = 0.1849929325344444
$ws.Range("R6").Value = 1.66493639281
$ws.Range("S6").Value = 0.07802833960396141
$ws.Range("T6").Value = 0.0780283396039614
$ws.Range("I7").Value = 0.6650661694281633
$ws.Range("J7").Value = 0.6650661694281633
$ws.Range("O7").Value = 0.03951584152489912
$ws.Range("P7").Value = 0.03951584152489912
$ws.Range("S7").Value = 0.02628064935469501
$ws.Range("T7").Value = 0.02628064935469501
$ws.Range("I8").Value = 0.6650661694281633
$ws.Range("J8").Value = 0.6650661694281633
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.9440163333333332
$ws.Range("N8").Value = 2.832049
$ws.Range("O8").Value = 0.8188575111383802
$ws.Range("P8").Value = 0.8188575111383801
$ws.Range("Q8").Value = 1.291147816721889
$ws.Range("R8").Value = 11.620330350497
$ws.Range("S8").Value = 0.5445944282402821
$ws.Range("T8").Value = 0.544594428240282
$ws.Range("I9").Value = 0.6650661694281633
$ws.Range("J9").Value = 0.6650661694281633
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.028017
$ws.Range("N9").Value = 0.084051
$ws.Range("O9").Value = 0.02430247240379386
$ws.Range("P9").Value = 0.02430247240379386
$ws.Range("Q9").Value = 0.038319345867
$ws.Range("R9").Value = 0.344874112803
$ws.Range("S9").Value = 0.01616275222922483
$ws.Range("T9").Value = 0.01616275222922483
$ws.Range("G10").Value = 0.5159453333333334
$ws.Range("H10").Value = 1.547836
$ws.Range("I10").Value = 0.2508834936018741
$ws.Range("J10").Value = 0.2508834936018741
$ws.Range("M10").Value = 0.1352566666666667
$ws.Range("N10").Value = 0.40577
$ws.Range("O10").Value = 0.1173241749329269
$ws.Range("P10").Value = 0.1173241749329268
$ws.Range("Q10").Value = 0.06978504596888889
$ws.Range("R10").Value = 0.6280654137199999
$ws.Range("S10").Value = 0.02943469889113012
$ws.Range("T10").Value = 0.02943469889113011
$ws.Range("G11").Value = 0.5159453333333334
$ws.Range("H11").Value = 1.547836
$ws.Range("I11").Value = 0.2508834936018741
$ws.Range("J11").Value = 0.2508834936018741
$ws.Range("O11").Value = 0.03951584152489912
$ws.Range("P11").Value = 0.03951584152489912
$ws.Range("Q11").Value = 0.02350423362355556
$ws.Range("R11").Value = 0.211538102612
$ws.Range("S11").Value = 0.009913872374384699
$ws.Range("T11").Value = 0.009913872374384699
$ws.Range("G12").Value = 0.5159453333333334
$ws.Range("H12").Value = 1.547836
$ws.Range("I12").Value = 0.2508834936018741
$ws.Range("J12").Value = 0.2508834936018741
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.9440163333333332
$ws.Range("N12").Value = 2.832049
$ws.Range("O12").Value = 0.8188575111383802
$ws.Range("P12").Value = 0.8188575111383801
$ws.Range("Q12").Value = 0.4870608217737777
$ws.Range("R12").Value = 4.383547395963999
$ws.Range("S12").Value = 0.2054378331565324
$ws.Range("T12").Value = 0.2054378331565324
$ws.Range("G13").Value = 0.5159453333333334
$ws.Range("H13").Value = 1.547836
$ws.Range("I13").Value = 0.2508834936018741
$ws.Range("J13").Value = 0.2508834936018741
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.028017
$ws.Range("N13").Value = 0.084051
$ws.Range("O13").Value = 0.02430247240379386
$ws.Range("P13").Value = 0.02430247240379386
$ws.Range("Q13").Value = 0.014455240404
$ws.Range("R13").Value = 0.130097163636
$ws.Range("S13").Value = 0.006097089179826941
$ws.Range("T13").Value = 0.00609708917982694
